# Rename the "YOLOv11" model labels to "YOLOv8" across the comparison tables
# on sheet "Hoja1". The three header rows (2, 25, 29, 33) each hold the same
# five model names in D:H; only the last three (F, G, H) need updating.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(2, 25, 29, 33)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "YOlOv8 n"
    $ws.Range("G$r").Value = "YOLOv8 pt"
    $ws.Range("H$r").Value = "YOLOv8 me"
}

# Restore the view/selection state captured at save time (best effort; the
# scroll position itself is not persisted by this runtime outside of
# freeze-pane contexts, but the active cell selection is).
$ws.Activate()
$ws.Range("I32").Select()
